# Applies the commit "Add files via upload":
# - Row 61 (2025-05-16 close): fills in the realized OHLCV + RV columns (C:H)
#   that were previously blank for that day.
# - Row 62 (2025-05-19): a brand-new trading day is appended with its full
#   gamma/OI wall snapshot (A, B, and I:II); C:H stay blank (no close yet).
# - Moves the active-cell selection on the frozen-pane view to G64.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 61: C61:H61 (Open/High/Low/Close/Volume/RV) ---
$ws.Range("C61").Value = 586.07
$ws.Range("D61").Value = 591.31
$ws.Range("E61").Value = 584.37
$ws.Range("F61").Value = 591.15
$ws.Range("G61").Value = 72389080
$ws.Range("H61").Value = 0.3774479169006175

# --- Row 62: new row -- dates (A62:B62) + full data block (I62:II62) ---
$ws.Range("A62").Value = 45804
$ws.Range("B62").Value = 45805
$ws.Range("I62").Value = 0.18960000000000002
$ws.Range("J62").Value = 4.6
$ws.Range("K62").Value = 600
$ws.Range("L62").Value = 274858800
$ws.Range("M62").Value = 4777
$ws.Range("N62").Value = 1224
$ws.Range("O62").Value = 6001
$ws.Range("P62").Value = 0.1020291606371589
$ws.Range("Q62").Value = 0.04144242150302037
$ws.Range("R62").Value = 45807
$ws.Range("S62").Value = 0.13201512290231018
$ws.Range("T62").Value = 45828
$ws.Range("U62").Value = 0.4060583541777592
$ws.Range("V62").Value = 45919
$ws.Range("W62").Value = 0.061317086811128435
$ws.Range("X62").Value = 47.333333333333336
$ws.Range("Y62").Value = 595
$ws.Range("Z62").Value = 140064190
$ws.Range("AA62").Value = -4308
$ws.Range("AB62").Value = 2442
$ws.Range("AC62").Value = 6750
$ws.Range("AD62").Value = 0.05199262945564612
$ws.Range("AE62").Value = 0.13715133500899862
$ws.Range("AF62").Value = 45805
$ws.Range("AG62").Value = 0.13715133500899862
$ws.Range("AH62").Value = 45807
$ws.Range("AI62").Value = 0.3433244451695464
$ws.Range("AJ62").Value = 45828
$ws.Range("AK62").Value = 0.2992694031334601
$ws.Range("AL62").Value = 9.333333333333334
$ws.Range("AM62").Value = 605
$ws.Range("AN62").Value = 136301055
$ws.Range("AO62").Value = -2746
$ws.Range("AP62").Value = 150
$ws.Range("AQ62").Value = 2896
$ws.Range("AR62").Value = 0.05059573219270851
$ws.Range("AS62").Value = 0.008162809346935355
$ws.Range("AT62").Value = 45807
$ws.Range("AU62").Value = 0.15008637576850495
$ws.Range("AV62").Value = 45814
$ws.Range("AW62").Value = 0.058440289009730736
$ws.Range("AX62").Value = 45828
$ws.Range("AY62").Value = 0.47862149858966124
$ws.Range("AZ62").Value = 12.333333333333334
$ws.Range("BA62").Value = 610
$ws.Range("BB62").Value = 93228130
$ws.Range("BC62").Value = 7349
$ws.Range("BD62").Value = 97
$ws.Range("BE62").Value = 7446
$ws.Range("BF62").Value = 0.03460681576020827
$ws.Range("BG62").Value = 0
$ws.Range("BH62").Value = 45807
$ws.Range("BI62").Value = 0.1153523324526665
$ws.Range("BJ62").Value = 45828
$ws.Range("BK62").Value = 0.26952348225803724
$ws.Range("BL62").Value = 45919
$ws.Range("BM62").Value = 0.11084572181819206
$ws.Range("BN62").Value = 47.333333333333336
$ws.Range("BO62").Value = 615
$ws.Range("BP62").Value = 84879225
$ws.Range("BQ62").Value = 1130
$ws.Range("BR62").Value = 93
$ws.Range("BS62").Value = 1223
$ws.Range("BT62").Value = 0.0315076544112197
$ws.Range("BU62").Value = 0
$ws.Range("BV62").Value = 45828
$ws.Range("BW62").Value = 0.17834499907269394
$ws.Range("BX62").Value = 45856
$ws.Range("BY62").Value = 0.26351291016121386
$ws.Range("BZ62").Value = 45919
$ws.Range("CA62").Value = 0.16017666897920774
$ws.Range("CB62").Value = 63.666666666666664
$ws.Range("CC62").Value = 575
$ws.Range("CD62").Value = -70907850
$ws.Range("CE62").Value = 0.03827602296737893
$ws.Range("CF62").Value = -874
$ws.Range("CG62").Value = -15296
$ws.Range("CH62").Value = 16170
$ws.Range("CI62").Value = 0.0367939955553692
$ws.Range("CJ62").Value = 45807
$ws.Range("CK62").Value = 0.1792632814793073
$ws.Range("CL62").Value = 45828
$ws.Range("CM62").Value = 0.3397731561071743
$ws.Range("CN62").Value = 45856
$ws.Range("CO62").Value = 0.09910268774372091
$ws.Range("CP62").Value = 26.333333333333332
$ws.Range("CQ62").Value = 585
$ws.Range("CR62").Value = -64832040
$ws.Range("CS62").Value = 0.034996303682343056
$ws.Range("CT62").Value = -6272
$ws.Range("CU62").Value = 29325
$ws.Range("CV62").Value = 35597
$ws.Range("CW62").Value = 0.17129981186403379
$ws.Range("CX62").Value = 45805
$ws.Range("CY62").Value = 0.17129981186403379
$ws.Range("CZ62").Value = 45807
$ws.Range("DA62").Value = 0.349653947600879
$ws.Range("DB62").Value = 45814
$ws.Range("DC62").Value = 0.1089576895539597
$ws.Range("DD62").Value = 4.666666666666667
$ws.Range("DE62").Value = 570
$ws.Range("DF62").Value = -59507430
$ws.Range("DG62").Value = 0.032122081792209094
$ws.Range("DH62").Value = -500
$ws.Range("DI62").Value = -1423
$ws.Range("DJ62").Value = 1923
$ws.Range("DK62").Value = 0
$ws.Range("DL62").Value = 45807
$ws.Range("DM62").Value = 0.307235423691558
$ws.Range("DN62").Value = 45828
$ws.Range("DO62").Value = 0.19314573268544544
$ws.Range("DP62").Value = 45838
$ws.Range("DQ62").Value = 0.1093147135715054
$ws.Range("DR62").Value = 20.333333333333332
$ws.Range("DS62").Value = 565
$ws.Range("DT62").Value = -49208110
$ws.Range("DU62").Value = 0.026562513861882832
$ws.Range("DV62").Value = -255
$ws.Range("DW62").Value = -3524
$ws.Range("DX62").Value = 3779
$ws.Range("DY62").Value = 0
$ws.Range("DZ62").Value = 45828
$ws.Range("EA62").Value = 0.17311127020611175
$ws.Range("EB62").Value = 45856
$ws.Range("EC62").Value = 0.4190361354987
$ws.Range("ED62").Value = 45919
$ws.Range("EE62").Value = 0.13096198048155544
$ws.Range("EF62").Value = 63.666666666666664
$ws.Range("EG62").Value = 555
$ws.Range("EH62").Value = -47440845
$ws.Range("EI62").Value = 0.025608545073808664
$ws.Range("EJ62").Value = -268
$ws.Range("EK62").Value = 5214
$ws.Range("EL62").Value = 5482
$ws.Range("EM62").Value = 0
$ws.Range("EN62").Value = 45828
$ws.Range("EO62").Value = 0.5315050522662549
$ws.Range("EP62").Value = 45856
$ws.Range("EQ62").Value = 0.2609327388513878
$ws.Range("ER62").Value = 45884
$ws.Range("ES62").Value = 0.06829900117069769
$ws.Range("ET62").Value = 52
$ws.Range("EU62").Value = 600
$ws.Range("EV62").Value = 351687600
$ws.Range("EW62").Value = 4777
$ws.Range("EX62").Value = 1224
$ws.Range("EY62").Value = 6001
$ws.Range("EZ62").Value = 0.07735410000012383
$ws.Range("FA62").Value = 313273200
$ws.Range("FB62").Value = 0.11628880591095067
$ws.Range("FC62").Value = 0.04144242150302037
$ws.Range("FD62").Value = 45807
$ws.Range("FE62").Value = 0.13201512290231018
$ws.Range("FF62").Value = 45828
$ws.Range("FG62").Value = 0.4060583541777592
$ws.Range("FH62").Value = 45919
$ws.Range("FI62").Value = 0.061317086811128435
$ws.Range("FJ62").Value = 47.333333333333336
$ws.Range("FK62").Value = -38414400
$ws.Range("FL62").Value = 0.020736074449839912
$ws.Range("FM62").Value = 0.005716606272647757
$ws.Range("FN62").Value = 45828
$ws.Range("FO62").Value = 0.1580969636386355
$ws.Range("FP62").Value = 45884
$ws.Range("FQ62").Value = 0.14564850680994626
$ws.Range("FR62").Value = 45919
$ws.Range("FS62").Value = 0.19990940897163564
$ws.Range("FT62").Value = 73
$ws.Range("FU62").Value = 590
$ws.Range("FV62").Value = 346823830
$ws.Range("FW62").Value = -2699
$ws.Range("FX62").Value = 19052
$ws.Range("FY62").Value = 21751
$ws.Range("FZ62").Value = 0.07628430808548821
$ws.Range("GA62").Value = 209620510
$ws.Range("GB62").Value = 0.07781233377877358
$ws.Range("GC62").Value = 0.15769697344978315
$ws.Range("GD62").Value = 45805
$ws.Range("GE62").Value = 0.15769697344978315
$ws.Range("GF62").Value = 45807
$ws.Range("GG62").Value = 0.31832114138067885
$ws.Range("GH62").Value = 45828
$ws.Range("GI62").Value = 0.21194013887286126
$ws.Range("GJ62").Value = 9.333333333333334
$ws.Range("GK62").Value = -137203320
$ws.Range("GL62").Value = 0.07406228545246599
$ws.Range("GM62").Value = 0.3207896864303284
$ws.Range("GN62").Value = 45805
$ws.Range("GO62").Value = 0.3207896864303284
$ws.Range("GP62").Value = 45807
$ws.Range("GQ62").Value = 0.16519600254571098
$ws.Range("GR62").Value = 45828
$ws.Range("GS62").Value = 0.10820131757744637
$ws.Range("GT62").Value = 9.333333333333334
$ws.Range("GU62").Value = 595
$ws.Range("GV62").Value = 326751390
$ws.Range("GW62").Value = -4308
$ws.Range("GX62").Value = 2442
$ws.Range("GY62").Value = 6750
$ws.Range("GZ62").Value = 0.0718693513710448
$ws.Range("HA62").Value = 233407790
$ws.Range("HB62").Value = 0.08664230834113461
$ws.Range("HC62").Value = 0.13715133500899862
$ws.Range("HD62").Value = 45805
$ws.Range("HE62").Value = 0.13715133500899862
$ws.Range("HF62").Value = 45807
$ws.Range("HG62").Value = 0.3433244451695464
$ws.Range("HH62").Value = 45828
$ws.Range("HI62").Value = 0.2992694031334601
$ws.Range("HJ62").Value = 9.333333333333334
$ws.Range("HK62").Value = -93343600
$ws.Range("HL62").Value = 0.05038682991315956
$ws.Range("HM62").Value = 0.029869964303926567
$ws.Range("HN62").Value = 45807
$ws.Range("HO62").Value = 0.049796022437531874
$ws.Range("HP62").Value = 45828
$ws.Range("HQ62").Value = 0.6324834268230495
$ws.Range("HR62").Value = 45884
$ws.Range("HS62").Value = 0.05813360530341662
$ws.Range("HT62").Value = 35.666666666666664
$ws.Range("HU62").Value = 589
$ws.Range("HV62").Value = 32103
$ws.Range("HW62").Value = 295314
$ws.Range("HX62").Value = 2693923955.5
$ws.Range("HY62").Value = -1852539645
$ws.Range("HZ62").Value = 841384310.5
$ws.Range("IA62").Value = 1.454178841878442
$ws.Range("IB62").Value = 4546463600.5
$ws.Range("IC62").Value = 0.13372276596102928
$ws.Range("ID62").Value = 45805
$ws.Range("IE62").Value = 0.13372276596102928
$ws.Range("IF62").Value = 45807
$ws.Range("IG62").Value = 0.1734910407977872
$ws.Range("IH62").Value = 45828
$ws.Range("II62").Value = 0.23821604617727324

# --- Restore the saved cursor position on the bottom-right (scrolling) pane ---
$ws.Range("G64").Select()

Write-Output "Applied SPY Walls row 61/62 update; selection -> G64"
